$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New supplier details (replaces "Distribuidora Isleña de Alimentos" with
# "EMSAR GLOBAL DE CENTRO AMERICA S.A.").
$ws.Range("A2").Value = "EMSAR GLOBAL DE CENTRO AMERICA S.A."

# Numeric-looking values must stay text (t="s"), so force the "Text"
# number format before assigning, then clear the format again so the
# cell keeps its original (unstyled) appearance.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3101775072"
$ws.Range("B2").ClearFormats()

$ws.Range("C2").Value = "EMSAR GLOBAL DE CENTRO AMERICA S.A."

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40301229"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "facelectronicaemsar@gmail.com"
